$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C numeric values (rows 3-14, excluding B-column label changes)
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 18
$ws.Range("C5").Value = 20
$ws.Range("C6").Value = 13
$ws.Range("C7").Value = 21
$ws.Range("C8").Value = 17

# Row 9: B column label + C column value
$ws.Range("B9").Value = "<number>"
$ws.Range("C9").Value = 11

# Row 10: B column label + C column value
$ws.Range("B10").Value = "<four>"
$ws.Range("C10").Value = 18

# Continue column C numeric values
$ws.Range("C11").Value = 15
$ws.Range("C12").Value = 19
$ws.Range("C13").Value = 18
$ws.Range("C14").Value = 21

# Row 15: B column label + C column value
$ws.Range("B15").Value = "<at>"
$ws.Range("C15").Value = 3
